# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
# - Shrink the "Status" column(s) width on each sheet to fit the new, shorter text

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$targetStoredWidth = 13.4101845877511

# This host stores column widths snapped to a 1/6-character pixel grid
# (stored = Round((ColumnWidth + 5/6) * 6) / 6, matching Excel's own
# digit-width rounding). Invert that so the ColumnWidth we send lands as
# close as possible to the width recorded in the target file.
$nearestStored = [Math]::Round($targetStoredWidth * 6) / 6
$newColumnWidth = $nearestStored - (5.0 / 6.0)

# Overview sheet: status shown in columns E (zh-cn) and F (de-de) of row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# Per-locale detail sheets: status shown in column C of row 2
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = $newColumnWidth
}
